# Applies the commit "modifies pullups and beautifies the report":
#  - fecha2 (C3) value updated from 21_09_2023 to 29_09_2023
#  - Introduccion (C6) text replaced with the new detailed paragraph
#  - Metodologia (C8) text replaced with the new short paragraph
#  - Columns B/C widened and rows 6-8 given explicit (taller) heights
#  - selection/active cell moved from C8 to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ---------------------------------------------
$ws.Range("C3").Value = "29_09_2023"

$nl = [char]10

$intro = "Evaluación de indicadores clave para el rendimiento en la escalada." + $nl +
    "Todas las evaluaciones de dedos son realizadas en una regleta de 20mm de profundidad con un radio de 10 mm" + $nl +
    "Fuerza máxima de dedos corresponde a la fuerza que es posible hacer luego de 5 segundo de fuerza isometrica maxima" + $nl +
    "Fuerza crítica es la fuerza que se puede ejercer despues de 4 minutos de 7 segundos de esfuerzo máximo y 3 segundos de descanso. Representa la fuerza petado" + $nl +
    "La potencia representa la fuerza que se puede ejercer por unidad de segundo. La usamos cuando caemos en una toma en movimiento dinámico y debemos ejercer mucha fuerza en poco tiempo" + $nl +
    "Todos los resultados se informan en % del peso corporal ya que así correlacionan mejor con el grado en la biografía" + $nl +
    "Comparamos los valores del deportista con referencias poblacionales para comparar con sus pares. Si los indicadores se encuentran por debajo de los valores de referencia, una posibilidad es entrenar para poder desarrollar y mantener mejores valores."

$ws.Range("C6").Value = $intro

$metodologia = "Medición con Tindeq de los indicadores clave informando el grado de escalada según escala IRCRA"

$ws.Range("C8").Value = $metodologia

# --- Column widths ------------------------------------------------------
# OOXML stored width = ColumnWidth + ~0.8333 (default padding), so subtract
# that back out to land on the target stored widths of 28.66 / 148.84
# (widths are pixel-quantized internally, so these land on the closest
# achievable stored value).
$ws.Columns.Item(2).ColumnWidth = 27.82
$ws.Columns.Item(3).ColumnWidth = 148.0

# --- Row heights ----------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 132.05
$ws.Rows.Item(7).RowHeight = 19.4
$ws.Rows.Item(8).RowHeight = 23.85

# --- View / selection state ------------------------------------------
$ws.Range("C9").Select()
